# Transaction report and update profile completed
# Update the generated Web registration ID on the "Reg ID" sheet (cell C2)
# to the new value produced by the latest registration run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reg ID")
$ws.Range("C2").Value = "Web00420180000001992"
